$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.142.44"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "3.666.39"
$ws.Range("E3").Value = "  +6.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.649"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("D8").Value = "3.660.39"
$ws.Range("E8").Value = "  +6.31%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.677"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.74%  "
$ws.Range("D15").Value = "4.251.30"
$ws.Range("E15").Value = "  +6.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.57%  "
$ws.Range("D17").Value = "3.659.23"
$ws.Range("E17").Value = "  +6.15%  "
$ws.Range("D18").Value = "71.166.93"
$ws.Range("E18").Value = "  +6.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.85%  "
$ws.Range("E32").Value = "  +9.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "628.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.42%  "
$ws.Range("D37").Value = "0.0₃0838"
$ws.Range("E37").Value = "  +10.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.414"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.26%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "3.341.78"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.15%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.94%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0457"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
